# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 930
$ws1.Range("F4").Value = 13311
$ws1.Range("F6").Value = 820
$ws1.Range("F7").Value = 4
$ws1.Range("F8").Value = 1377
$ws1.Range("F9").Value = 120
$ws1.Range("F13").Value = 27
$ws1.Range("F14").Value = 13298
$ws1.Range("F17").Value = 8877
$ws1.Range("F19").Value = 7956
$ws1.Range("F20").Value = 240
$ws1.Range("F26").Value = 15
$ws1.Range("F27").Value = 1010
$ws1.Range("F32").Value = 141
$ws1.Range("F33").Value = 365

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 930
$ws4.Range("F5").Value = 13311
$ws4.Range("F7").Value = 822
$ws4.Range("F8").Value = 4
$ws4.Range("F9").Value = 1377
$ws4.Range("F10").Value = 120
$ws4.Range("F14").Value = 27
$ws4.Range("F15").Value = 13298
$ws4.Range("F18").Value = 8877
$ws4.Range("F20").Value = 7956
$ws4.Range("F21").Value = 240
$ws4.Range("F27").Value = 15
$ws4.Range("F28").Value = 1010
$ws4.Range("F35").Value = 141
$ws4.Range("F36").Value = 365
